$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 925.5
$ws.Range("J41").Value = 925.5
$ws.Range("L41").Value = 925.5
$ws.Range("N41").Value = -1805.5

$ws.Range("H98").Value = 83334696
$ws.Range("I98").Value = 89287000
$ws.Range("J98").Value = 2500
$ws.Range("K98").Value = 89287000
$ws.Range("L98").Value = 2500
$ws.Range("M98").Value = -89285502
$ws.Range("N98").Value = -5496

$ws.Range("H122").Value = 83334696
$ws.Range("I122").Value = 89287000
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 267861000
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -267858550
$ws.Range("N122").Value = -12400

$ws.Range("H135").Value = 5209233.5
$ws.Range("I135").Value = 6250766.5
$ws.Range("J135").Value = 1568.5
$ws.Range("K135").Value = 56256898.5
$ws.Range("L135").Value = 14116.5
$ws.Range("M135").Value = -56254363.5
$ws.Range("N135").Value = -19186.5

$ws.Range("H137").Value = 1623.6364
$ws.Range("I137").Value = 1200.4
$ws.Range("J137").Value = 1976.3334
$ws.Range("K137").Value = 3601.2
$ws.Range("L137").Value = 5929.0002
$ws.Range("M137").Value = -1051.2
$ws.Range("N137").Value = -11029.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H42").Value = 9080
$ws.Range("J42").Value = 9080
$ws.Range("L42").Value = 9080
$ws.Range("N42").Value = -10052

$ws.Range("H61").Value = 4809215
$ws.Range("I61").Value = 5683383
$ws.Range("J61").Value = 1291
$ws.Range("K61").Value = 5683383
$ws.Range("L61").Value = 1291
$ws.Range("M61").Value = -5683171
$ws.Range("N61").Value = -1715

$ws.Range("H107").Value = 30228
$ws.Range("J107").Value = 30228
$ws.Range("L107").Value = 30228
$ws.Range("N107").Value = -37908

$ws.Range("H132").Value = 2559125.5
$ws.Range("I132").Value = 1344.7059
$ws.Range("J132").Value = 9806171
$ws.Range("K132").Value = 4034.1177
$ws.Range("L132").Value = 29418513
$ws.Range("M132").Value = -1504.1177
$ws.Range("N132").Value = -29423573

$ws.Range("H136").Value = 4809215
$ws.Range("I136").Value = 5683383
$ws.Range("J136").Value = 1291
$ws.Range("K136").Value = 17050149
$ws.Range("L136").Value = 3873
$ws.Range("M136").Value = -17047599
$ws.Range("N136").Value = -8973

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5057036.5
$ws.Range("I134").Value = 2211.7144
$ws.Range("J134").Value = 13902980
$ws.Range("K134").Value = 6635.1432
$ws.Range("L134").Value = 41708940
$ws.Range("M134").Value = -4100.1432
$ws.Range("N134").Value = -41714010

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 71479816
$ws.Range("I86").Value = 136412540
$ws.Range("J86").Value = 53810
$ws.Range("K86").Value = 136412540
$ws.Range("L86").Value = 53810
$ws.Range("M86").Value = -136411417
$ws.Range("N86").Value = -56056

$ws.Range("H89").Value = 71479816
$ws.Range("I89").Value = 136412540
$ws.Range("J89").Value = 53810
$ws.Range("K89").Value = 682062700
$ws.Range("L89").Value = 269050
$ws.Range("M89").Value = -682057084
$ws.Range("N89").Value = -280282

$ws.Range("H103").Value = 3934.1428
$ws.Range("I103").Value = 3934.1428
$ws.Range("K103").Value = 3934.1428
$ws.Range("M103").Value = -2762.1428

$ws.Range("H122").Value = 13890345
$ws.Range("I122").Value = 20834976
$ws.Range("J122").Value = 1084
$ws.Range("K122").Value = 62504928
$ws.Range("L122").Value = 3252
$ws.Range("M122").Value = -62502478
$ws.Range("N122").Value = -8152

$ws.Range("H132").Value = 19610694
$ws.Range("I132").Value = 2838.3333
$ws.Range("K132").Value = 8514.999899999999
$ws.Range("M132").Value = -5984.999899999999

$ws.Range("H134").Value = 16667793
$ws.Range("I134").Value = 1060.6957
$ws.Range("J134").Value = 71429910
$ws.Range("K134").Value = 3182.0871
$ws.Range("L134").Value = 214289730
$ws.Range("M134").Value = -647.0870999999997
$ws.Range("N134").Value = -214294800

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3365.5
$ws.Range("I131").Value = 2546
$ws.Range("J131").Value = 3410.5276
$ws.Range("K131").Value = 7638
$ws.Range("L131").Value = 10231.5828
$ws.Range("M131").Value = -2598
$ws.Range("N131").Value = -20311.5828

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9658.105
$ws.Range("I70").Value = 11501.714
$ws.Range("J70").Value = 4496
$ws.Range("K70").Value = 11501.714
$ws.Range("L70").Value = 4496
$ws.Range("M70").Value = -11231.714
$ws.Range("N70").Value = -5036

$ws.Range("H73").Value = 9658.105
$ws.Range("I73").Value = 11501.714
$ws.Range("J73").Value = 4496
$ws.Range("K73").Value = 11501.714
$ws.Range("L73").Value = 4496
$ws.Range("M73").Value = -10565.714
$ws.Range("N73").Value = -6368

$ws.Range("H132").Value = 7247.8945
$ws.Range("I132").Value = 1700.9231
$ws.Range("J132").Value = 19266.334
$ws.Range("K132").Value = 5102.7693
$ws.Range("L132").Value = 57799.00199999999
$ws.Range("M132").Value = -2572.7693
$ws.Range("N132").Value = -62859.00199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 49800
$ws.Range("J6").Value = 49800
$ws.Range("L6").Value = 49800
$ws.Range("N6").Value = -50024

$ws.Range("H132").Value = 42337756
$ws.Range("I132").Value = 81633990
$ws.Range("J132").Value = 18731
$ws.Range("K132").Value = 244901970
$ws.Range("L132").Value = 56193
$ws.Range("M132").Value = -244899440
$ws.Range("N132").Value = -61253

$ws.Range("H136").Value = 107663430
$ws.Range("I136").Value = 113556160
$ws.Range("J136").Value = 100002890
$ws.Range("K136").Value = 340668480
$ws.Range("L136").Value = 300008670
$ws.Range("M136").Value = -340665930
$ws.Range("N136").Value = -300013770

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 26034.623
$ws.Range("I132").Value = 38317.035
$ws.Range("J132").Value = 7611
$ws.Range("K132").Value = 114951.105
$ws.Range("L132").Value = 22833
$ws.Range("M132").Value = -112421.105
$ws.Range("N132").Value = -27893

$ws.Range("H136").Value = 9617634
$ws.Range("J136").Value = 1939.3103
$ws.Range("L136").Value = 5817.9309
$ws.Range("N136").Value = -10917.9309
